$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.545.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.74%  "

$ws.Range("D3").Value = "'3.745.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.11%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'610.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.01%  "

$ws.Range("D6").Value = "'176.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.55%  "

$ws.Range("D7").Value = "'3.738.68"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.12%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").Value = "'0.536"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.19%  "

$ws.Range("E10").Value = "  +5.49%  "

$ws.Range("D11").Value = "'6.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.13%  "

$ws.Range("D12").Value = "'0.492"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.89%  "

$ws.Range("D13").Value = "'40.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.36%  "

$ws.Range("D14").Value = "'0.0000252"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.48%  "

$ws.Range("D15").Value = "'4.374.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.36%  "

$ws.Range("D16").Value = "'3.750.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.91%  "

$ws.Range("D17").Value = "'69.654.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.48%  "

$ws.Range("D18").Value = "'0.123"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.51%  "

$ws.Range("D19").Value = "'7.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.15%  "

$ws.Range("D20").Value = "'512.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.05%  "

$ws.Range("D21").Value = "'16.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.59%  "

$ws.Range("D22").Value = "'9.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.08%  "

$ws.Range("D23").Value = "'0.723"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.40%  "

$ws.Range("D24").Value = "'87.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.81%  "

$ws.Range("D25").Value = "'2.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.12%  "

$ws.Range("D26").Value = "'13.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("D27").Value = "'11.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.10%  "

$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("D29").Value = "'0.0000126"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.90%  "

$ws.Range("D30").Value = "'2.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.97%  "

$ws.Range("D31").Value = "'2.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.40%  "

$ws.Range("D32").Value = "'7.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.10%  "

$ws.Range("D33").Value = "'31.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("D34").Value = "'0.114"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.01%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("D36").Value = "'6.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.28%  "

$ws.Range("E37").Value = "  +2.86%  "

$ws.Range("E38").Value = "  +2.30%  "

$ws.Range("E39").Value = "  +4.21%  "

$ws.Range("D40").Value = "'0.132"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.41%  "

$ws.Range("D41").Value = "'51.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.08%  "

$ws.Range("E42").Value = "  -4.23%  "

$ws.Range("D43").Value = "'8.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.77%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "'3.071.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.20%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'418.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.09%  "

$ws.Range("D46").Value = "'2.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.38%  "

$ws.Range("D47").Value = "'0.0362"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.41%  "

$ws.Range("D48").Value = "'27.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.40%  "

$ws.Range("E49").Value = "  +2.45%  "

$ws.Range("D50").Value = "'135.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.24%  "

$ws.Range("E51").Value = "  -0.06%  "
